$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("B1").Value = "C/A"
$ws.Range("C1").Value = "FFR"

# New D1 header ("LF") needs the same bold/border/centered style as B1/C1.
# Copy+PasteSpecial(formats) from B1 reuses the existing style index instead
# of accumulating a new one.
$ws.Range("D1").Value = "LF"
$ws.Range("B1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null

# --- Row labels (column A) ---
$ws.Range("A2").Value = "C/A Lag"
$ws.Range("A3").Value = "FFR Lag"
$ws.Range("A4").Value = "LF Lag"
$ws.Range("A5").Value = "r2"

# --- Row 2 data ---
$ws.Range("B2").Value = "-0.699***"
$ws.Range("C2").Value = "-0.029***"
$ws.Range("D2").Value = "-0.235***"

# --- Row 3 data ---
$ws.Range("B3").Value = "9.766***"
$ws.Range("C3").Value = "0.395***"
$ws.Range("D3").Value = "10.288***"

# --- Row 4 data ---
# "0.041" and "-0.006" must stay text (shared string), not become numbers.
# Stage each value as text in a scratch cell, then copy only the VALUE
# (xlPasteValues) into the destination so the destination cell never
# directly goes through the auto-number-conversion / quote-prefix style
# path itself.
$ws.Range("Z1").Formula = "'0.041"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4163) | Out-Null

$ws.Range("Z1").Formula = "'-0.006"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("C4").PasteSpecial(-4163) | Out-Null

$ws.Range("Z1").Clear() | Out-Null

$ws.Range("D4").Value = "-0.136***"

# --- Row 5 data (numeric) ---
$ws.Range("B5").Value = 0.7423139539042245
$ws.Range("C5").Value = 0.462781965776386
$ws.Range("D5").Value = 0.8046210953053092
